$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 14

$ws.Cells.Item($row, 1).Value = 42620.886423611111

$ws.Cells.Item($row, 2).Value = -6
$ws.Cells.Item($row, 3).Value = 56
$ws.Cells.Item($row, 4).Value = 40
$ws.Cells.Item($row, 5).Value = 56
$ws.Cells.Item($row, 6).Value = 70
$ws.Cells.Item($row, 7).Value = 16990
$ws.Cells.Item($row, 8).Value = 11886
$ws.Cells.Item($row, 9).Value = 1835
$ws.Cells.Item($row, 10).Value = 218
$ws.Cells.Item($row, 11).Value = 155
$ws.Cells.Item($row, 12).Value = 5
$ws.Cells.Item($row, 13).Value = 12
$ws.Cells.Item($row, 14).Value = "Noun"
